$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price/volume columns for rows with changed crypto data.
# Some new price strings (e.g. "319.12") parse as plain numbers, unlike
# the original text values ("44.267.94" has two dots so it always stays
# text) -- force NumberFormat to Text first so the cell keeps the same
# text type as before the edit.
$ws.Range("D2").Value = "44.138.91"
$ws.Range("E2").Value = "  +3.53%  "
$ws.Range("D3").Value = "2.272.20"
$ws.Range("E3").Value = "  +2.39%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.12"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.82"
$ws.Range("E6").Value = "  +5.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.589"
$ws.Range("E7").Value = "  +2.45%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.576"
$ws.Range("E9").Value = "  +2.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.78"
$ws.Range("E10").Value = "  +6.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0843"
$ws.Range("E11").Value = "  +2.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.90"
$ws.Range("E12").Value = "  +2.94%  "
$ws.Range("E13").Value = "  +2.52%  "
$ws.Range("D14").Value = "2.622.56"
$ws.Range("E14").Value = "  +2.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.881"
$ws.Range("E15").Value = "  +2.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.63"
$ws.Range("E16").Value = "  +4.52%  "
$ws.Range("D17").Value = "2.283.31"
$ws.Range("E17").Value = "  +2.75%  "
$ws.Range("D18").Value = "44.114.23"
$ws.Range("E18").Value = "  +3.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.61"
$ws.Range("E19").Value = "  +6.32%  "
$ws.Range("D20").Value = "0.0₃0999"
$ws.Range("E20").Value = "  +3.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.67"
$ws.Range("E21").Value = "  +2.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.34"
$ws.Range("E22").Value = "  +1.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.22"
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "239.54"
$ws.Range("E24").Value = "  +1.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.23"
$ws.Range("E25").Value = "  +3.92%  "
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.06"
$ws.Range("E27").Value = "  +1.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "40.10"
$ws.Range("E28").Value = "  +18.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.25"
$ws.Range("E29").Value = "  +1.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.20"
$ws.Range("E30").Value = "  +0.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.52"
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0890"
$ws.Range("E32").Value = "  +0.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.57"
$ws.Range("E33").Value = "  +0.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "161.12"
$ws.Range("E34").Value = "  +2.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.38"
$ws.Range("E35").Value = "  +5.42%  "
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.55"
$ws.Range("E39").Value = "  +2.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.95"
$ws.Range("E40").Value = "  +10.91%  "
$ws.Range("E41").Value = "  +5.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.12"
$ws.Range("E42").Value = "  +34.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0329"
$ws.Range("E43").Value = "  +2.30%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "1.793.56"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("E46").Value = "  +1.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.53"
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.42"
$ws.Range("E48").Value = "  +1.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.31"
$ws.Range("E49").Value = "  -1.82%  "

# Row 50 and 51 swapped: MultiversX (row50) and FraxShare (row51) traded places
# with updated price/volume data
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.84"
$ws.Range("E50").Value = "  +4.85%  "

$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "59.91"
$ws.Range("E51").Value = "  -0.80%  "
